$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.223.94'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.325.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("E7").Value = '  -1.41%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.61'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.60'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.40%  '
$ws.Range("E12").Value = '  +0.87%  '
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.685.67'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.326.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.788'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.014.76'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("E19").Value = '  -1.40%  '
$ws.Range("E20").Value = '  -0.51%  '
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.88'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.74'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.23'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.80%  '
$ws.Range("E28").Value = '  +1.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '165.39'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.29'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.01'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.54'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0698'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.31%  '
$ws.Range("E37").Value = '  -1.80%  '
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("E39").Value = '  +1.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.110'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.990.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.70'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.08%  '
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.20%  '
$ws.Range("E46").Value = '  -3.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.78'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("E48").Value = '  -2.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.552.20'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.09'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.31%  '
